$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 130.92308
$ws.Range("I9").Value = 138.4
$ws.Range("K9").Value = 138.4
$ws.Range("M9").Value = 30.59999999999999
$ws.Range("H18").Value = 5856.357
$ws.Range("I18").Value = 1495
$ws.Range("J18").Value = 6583.25
$ws.Range("K18").Value = 1495
$ws.Range("L18").Value = 6583.25
$ws.Range("M18").Value = -1211
$ws.Range("N18").Value = -7151.25
$ws.Range("H33").Value = 212.11765
$ws.Range("I33").Value = 213.8
$ws.Range("K33").Value = 213.8
$ws.Range("M33").Value = 15.19999999999999
$ws.Range("H82").Value = 7933
$ws.Range("I82").Value = 7899.5
$ws.Range("K82").Value = 23698.5
$ws.Range("M82").Value = -23292.5
$ws.Range("H85").Value = 7933
$ws.Range("I85").Value = 7899.5
$ws.Range("K85").Value = 23698.5
$ws.Range("M85").Value = -22294.5
$ws.Range("H106").Value = 3241.9473
$ws.Range("I106").Value = 3023.4119
$ws.Range("K106").Value = 3023.4119
$ws.Range("M106").Value = -2392.4119
$ws.Range("H137").Value = 2054.923
$ws.Range("I137").Value = 1059.75
$ws.Range("J137").Value = 2497.2222
$ws.Range("K137").Value = 3179.25
$ws.Range("L137").Value = 7491.6666
$ws.Range("M137").Value = -629.25
$ws.Range("N137").Value = -12591.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3532.2068
$ws.Range("I2").Value = 2625.8262
$ws.Range("J2").Value = 7006.6665
$ws.Range("K2").Value = 2625.8262
$ws.Range("L2").Value = 7006.6665
$ws.Range("M2").Value = -2512.8262
$ws.Range("N2").Value = -7232.6665
$ws.Range("H45").Value = 4091.4167
$ws.Range("I45").Value = 4210.8
$ws.Range("K45").Value = 4210.8
$ws.Range("M45").Value = -3833.8
$ws.Range("H61").Value = 1963176.6
$ws.Range("I61").Value = 2224000.2
$ws.Range("K61").Value = 2224000.2
$ws.Range("M61").Value = -2223788.2
$ws.Range("H74").Value = 4633955
$ws.Range("I74").Value = 6251738
$ws.Range("K74").Value = 6251738
$ws.Range("M74").Value = -6250864
$ws.Range("H77").Value = 4633955
$ws.Range("I77").Value = 6251738
$ws.Range("K77").Value = 31258690
$ws.Range("M77").Value = -31254322
$ws.Range("H116").Value = 3532.2068
$ws.Range("I116").Value = 2625.8262
$ws.Range("J116").Value = 7006.6665
$ws.Range("K116").Value = 2625.8262
$ws.Range("L116").Value = 7006.6665
$ws.Range("M116").Value = -331.8262
$ws.Range("N116").Value = -11594.6665
$ws.Range("H136").Value = 1963176.6
$ws.Range("I136").Value = 2224000.2
$ws.Range("K136").Value = 6672000.600000001
$ws.Range("M136").Value = -6669450.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3532.2068
$ws.Range("I3").Value = 2625.8262
$ws.Range("J3").Value = 7006.6665
$ws.Range("K3").Value = 2625.8262
$ws.Range("L3").Value = 7006.6665
$ws.Range("M3").Value = -2511.8262
$ws.Range("N3").Value = -7234.6665
$ws.Range("H20").Value = 1269
$ws.Range("I20").Value = 924.8182
$ws.Range("J20").Value = 1809.8572
$ws.Range("K20").Value = 924.8182
$ws.Range("L20").Value = 1809.8572
$ws.Range("M20").Value = -677.8182
$ws.Range("N20").Value = -2303.8572
$ws.Range("H62").Value = 90113.125
$ws.Range("I62").Value = 73333.336
$ws.Range("K62").Value = 73333.336
$ws.Range("M62").Value = -72647.336
$ws.Range("H65").Value = 90113.125
$ws.Range("I65").Value = 73333.336
$ws.Range("K65").Value = 220000.008
$ws.Range("M65").Value = -216568.008
$ws.Range("H86").Value = 2153.4375
$ws.Range("I86").Value = 1995.2222
$ws.Range("K86").Value = 1995.2222
$ws.Range("M86").Value = -872.2221999999999
$ws.Range("H89").Value = 2153.4375
$ws.Range("I89").Value = 1995.2222
$ws.Range("K89").Value = 9976.110999999999
$ws.Range("M89").Value = -4360.110999999999
$ws.Range("H99").Value = 8367.053
$ws.Range("I99").Value = 3982.1538
$ws.Range("K99").Value = 3982.1538
$ws.Range("M99").Value = -2484.1538
$ws.Range("H134").Value = 431253.22
$ws.Range("I134").Value = 474863.44
$ws.Range("J134").Value = 292245.7
$ws.Range("K134").Value = 1424590.32
$ws.Range("L134").Value = 876737.1000000001
$ws.Range("M134").Value = -1422055.32
$ws.Range("N134").Value = -881807.1000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 866.06665
$ws.Range("I22").Value = 874.3333
$ws.Range("J22").Value = 833
$ws.Range("K22").Value = 874.3333
$ws.Range("L22").Value = 833
$ws.Range("M22").Value = -524.3333
$ws.Range("N22").Value = -1533
$ws.Range("H105").Value = 17545.762
$ws.Range("I105").Value = 20075.666
$ws.Range("J105").Value = 2366.3333
$ws.Range("K105").Value = 20075.666
$ws.Range("L105").Value = 2366.3333
$ws.Range("M105").Value = -18328.666
$ws.Range("N105").Value = -5860.3333
$ws.Range("H107").Value = 1174.8462
$ws.Range("I107").Value = 1077.8
$ws.Range("K107").Value = 1077.8
$ws.Range("M107").Value = 842.2
$ws.Range("H134").Value = 18446.2
$ws.Range("I134").Value = 20022.777
$ws.Range("K134").Value = 60068.33099999999
$ws.Range("M134").Value = -57533.33099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1339
$ws.Range("J132").Value = 1403.75
$ws.Range("L132").Value = 12633.75
$ws.Range("N132").Value = -17693.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 25000
$ws.Range("J34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("N34").Value = -25536
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25630
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27184
$ws.Range("H96").Value = 42252
$ws.Range("J96").Value = 42252
$ws.Range("L96").Value = 42252
$ws.Range("N96").Value = -47744
$ws.Range("H113").Value = 2076.0833
$ws.Range("J113").Value = 2954
$ws.Range("L113").Value = 2954
$ws.Range("N113").Value = -7294

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5889.448
$ws.Range("I100").Value = 1947.619
$ws.Range("J100").Value = 16236.75
$ws.Range("K100").Value = 1947.619
$ws.Range("L100").Value = 16236.75
$ws.Range("M100").Value = -1406.619
$ws.Range("N100").Value = -17318.75
$ws.Range("H136").Value = 44819.766
$ws.Range("I136").Value = 1366.7059
$ws.Range("J136").Value = 101643
$ws.Range("K136").Value = 4100.1177
$ws.Range("L136").Value = 304929
$ws.Range("M136").Value = -1550.1177
$ws.Range("N136").Value = -310029

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3825.5715
$ws.Range("I62").Value = 3595.8
$ws.Range("J62").Value = 4400
$ws.Range("K62").Value = 3595.8
$ws.Range("L62").Value = 4400
$ws.Range("M62").Value = -2971.8
$ws.Range("N62").Value = -5648
$ws.Range("H65").Value = 3825.5715
$ws.Range("I65").Value = 3595.8
$ws.Range("J65").Value = 4400
$ws.Range("K65").Value = 17979
$ws.Range("L65").Value = 22000
$ws.Range("M65").Value = -14859
$ws.Range("N65").Value = -28240
$ws.Range("H122").Value = 2858.2341
$ws.Range("I122").Value = 2518.0256
$ws.Range("K122").Value = 7554.0768
$ws.Range("M122").Value = -5104.0768
$ws.Range("H132").Value = 3032553
$ws.Range("J132").Value = 9335.666999999999
$ws.Range("L132").Value = 28007.001
$ws.Range("N132").Value = -33067.001
